$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(70, 8).Value = 101500.5  # H70
$ws.Cells.Item(70, 9).Value = 101500.5  # I70
$ws.Cells.Item(70, 10).Value = 0  # J70
$ws.Cells.Item(70, 11).Value = 304501.5  # K70
$ws.Cells.Item(70, 12).Value = 0  # L70
$ws.Cells.Item(70, 13).Value = -304231.5  # M70
$ws.Cells.Item(70, 14).ClearContents()  # N70

$ws.Cells.Item(73, 8).Value = 101500.5  # H73
$ws.Cells.Item(73, 9).Value = 101500.5  # I73
$ws.Cells.Item(73, 10).Value = 0  # J73
$ws.Cells.Item(73, 11).Value = 304501.5  # K73
$ws.Cells.Item(73, 12).Value = 0  # L73
$ws.Cells.Item(73, 13).Value = -303565.5  # M73
$ws.Cells.Item(73, 14).ClearContents()  # N73

$ws.Cells.Item(86, 8).Value = 3407.9473  # H86
$ws.Cells.Item(86, 9).Value = 2785.6  # I86
$ws.Cells.Item(86, 10).Value = 4099.4443  # J86
$ws.Cells.Item(86, 11).Value = 2785.6  # K86
$ws.Cells.Item(86, 12).Value = 4099.4443  # L86
$ws.Cells.Item(86, 13).Value = -1662.6  # M86
$ws.Cells.Item(86, 14).Value = -6345.4443  # N86

$ws.Cells.Item(89, 8).Value = 3407.9473  # H89
$ws.Cells.Item(89, 9).Value = 2785.6  # I89
$ws.Cells.Item(89, 10).Value = 4099.4443  # J89
$ws.Cells.Item(89, 11).Value = 13928  # K89
$ws.Cells.Item(89, 12).Value = 20497.2215  # L89
$ws.Cells.Item(89, 13).Value = -8312  # M89
$ws.Cells.Item(89, 14).Value = -31729.2215  # N89

$ws.Cells.Item(100, 8).Value = 32035.758  # H100
$ws.Cells.Item(100, 9).Value = 36126.965  # I100
$ws.Cells.Item(100, 10).Value = 2374.5  # J100
$ws.Cells.Item(100, 11).Value = 36126.965  # K100
$ws.Cells.Item(100, 12).Value = 2374.5  # L100
$ws.Cells.Item(100, 13).Value = -35585.965  # M100
$ws.Cells.Item(100, 14).Value = -3456.5  # N100

$ws.Cells.Item(103, 8).Value = 1083  # H103
$ws.Cells.Item(103, 9).Value = 1106  # I103
$ws.Cells.Item(103, 10).Value = 1048.5  # J103
$ws.Cells.Item(103, 11).Value = 3318  # K103
$ws.Cells.Item(103, 12).Value = 3145.5  # L103
$ws.Cells.Item(103, 13).Value = -2732  # M103
$ws.Cells.Item(103, 14).Value = -4317.5  # N103

$ws.Cells.Item(113, 8).Value = 4125.75  # H113
$ws.Cells.Item(113, 9).Value = 4167.6665  # I113
$ws.Cells.Item(113, 10).Value = 4000  # J113
$ws.Cells.Item(113, 11).Value = 4167.6665  # K113
$ws.Cells.Item(113, 12).Value = 4000  # L113
$ws.Cells.Item(113, 13).Value = -913.6665000000003  # M113
$ws.Cells.Item(113, 14).Value = -10508  # N113

$ws.Cells.Item(137, 8).Value = 9108.536  # H137
$ws.Cells.Item(137, 9).Value = 3887.7666  # I137
$ws.Cells.Item(137, 10).Value = 15132.5  # J137
$ws.Cells.Item(137, 11).Value = 11663.2998  # K137
$ws.Cells.Item(137, 12).Value = 45397.5  # L137
$ws.Cells.Item(137, 13).Value = -9113.2998  # M137
$ws.Cells.Item(137, 14).Value = -50497.5  # N137

$ws.Cells.Item(138, 8).Value = 2993.543  # H138
$ws.Cells.Item(138, 9).Value = 1822.0416  # I138
$ws.Cells.Item(138, 10).Value = 5549.5454  # J138
$ws.Cells.Item(138, 11).Value = 5466.1248  # K138
$ws.Cells.Item(138, 12).Value = 16648.6362  # L138
$ws.Cells.Item(138, 13).Value = -326.1247999999996  # M138
$ws.Cells.Item(138, 14).Value = -26928.6362  # N138

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(32, 8).Value = 3586.907  # H32
$ws.Cells.Item(32, 9).Value = 1964.0883  # I32
$ws.Cells.Item(32, 10).Value = 9717.556  # J32
$ws.Cells.Item(32, 11).Value = 1964.0883  # K32
$ws.Cells.Item(32, 12).Value = 9717.556  # L32
$ws.Cells.Item(32, 13).Value = -1677.0883  # M32
$ws.Cells.Item(32, 14).Value = -10291.556  # N32

$ws.Cells.Item(61, 8).Value = 7971.4546  # H61
$ws.Cells.Item(61, 9).Value = 4884.5713  # I61
$ws.Cells.Item(61, 10).Value = 13373.5  # J61
$ws.Cells.Item(61, 11).Value = 4884.5713  # K61
$ws.Cells.Item(61, 12).Value = 13373.5  # L61
$ws.Cells.Item(61, 13).Value = -4672.5713  # M61
$ws.Cells.Item(61, 14).Value = -13797.5  # N61

$ws.Cells.Item(74, 8).Value = 2986.0527  # H74
$ws.Cells.Item(74, 9).Value = 1570.8334  # I74
$ws.Cells.Item(74, 10).Value = 5412.143  # J74
$ws.Cells.Item(74, 11).Value = 1570.8334  # K74
$ws.Cells.Item(74, 12).Value = 5412.143  # L74
$ws.Cells.Item(74, 13).Value = -696.8334  # M74
$ws.Cells.Item(74, 14).Value = -7160.143  # N74

$ws.Cells.Item(77, 8).Value = 2986.0527  # H77
$ws.Cells.Item(77, 9).Value = 1570.8334  # I77
$ws.Cells.Item(77, 10).Value = 5412.143  # J77
$ws.Cells.Item(77, 11).Value = 7854.166999999999  # K77
$ws.Cells.Item(77, 12).Value = 27060.715  # L77
$ws.Cells.Item(77, 13).Value = -3486.166999999999  # M77
$ws.Cells.Item(77, 14).Value = -35796.715  # N77

$ws.Cells.Item(122, 8).Value = 4438.6665  # H122
$ws.Cells.Item(122, 9).Value = 4627.857  # I122
$ws.Cells.Item(122, 10).Value = 3776.5  # J122
$ws.Cells.Item(122, 11).Value = 13883.571  # K122
$ws.Cells.Item(122, 12).Value = 11329.5  # L122
$ws.Cells.Item(122, 13).Value = -11433.571  # M122
$ws.Cells.Item(122, 14).Value = -16229.5  # N122

$ws.Cells.Item(132, 8).Value = 9281.695  # H132
$ws.Cells.Item(132, 9).Value = 10712.158  # I132
$ws.Cells.Item(132, 10).Value = 2487  # J132
$ws.Cells.Item(132, 11).Value = 32136.474  # K132
$ws.Cells.Item(132, 12).Value = 7461  # L132
$ws.Cells.Item(132, 13).Value = -29606.474  # M132
$ws.Cells.Item(132, 14).Value = -12521  # N132

$ws.Cells.Item(136, 8).Value = 7971.4546  # H136
$ws.Cells.Item(136, 9).Value = 4884.5713  # I136
$ws.Cells.Item(136, 10).Value = 13373.5  # J136
$ws.Cells.Item(136, 11).Value = 14653.7139  # K136
$ws.Cells.Item(136, 12).Value = 40120.5  # L136
$ws.Cells.Item(136, 13).Value = -12103.7139  # M136
$ws.Cells.Item(136, 14).Value = -45220.5  # N136

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Cells.Item(107, 8).Value = 1292.8182  # H107
$ws.Cells.Item(107, 9).Value = 1560.2858  # I107
$ws.Cells.Item(107, 10).Value = 824.75  # J107
$ws.Cells.Item(107, 11).Value = 1560.2858  # K107
$ws.Cells.Item(107, 12).Value = 824.75  # L107
$ws.Cells.Item(107, 13).Value = 359.7141999999999  # M107
$ws.Cells.Item(107, 14).Value = -4664.75  # N107

$ws.Cells.Item(132, 8).Value = 79757.27  # H132
$ws.Cells.Item(132, 9).Value = 0  # I132
$ws.Cells.Item(132, 10).Value = 79757.27  # J132
$ws.Cells.Item(132, 11).Value = 0  # K132
$ws.Cells.Item(132, 12).Value = 79757.27  # L132
$ws.Cells.Item(132, 14).Value = -89877.27  # N132

$ws.Cells.Item(133, 8).Value = 75000  # H133
$ws.Cells.Item(133, 9).Value = 0  # I133
$ws.Cells.Item(133, 10).Value = 75000  # J133
$ws.Cells.Item(133, 11).Value = 0  # K133
$ws.Cells.Item(133, 12).Value = 75000  # L133
$ws.Cells.Item(133, 14).Value = -85120  # N133

$ws.Cells.Item(134, 8).Value = 6892.4062  # H134
$ws.Cells.Item(134, 9).Value = 4524.174  # I134
$ws.Cells.Item(134, 10).Value = 12944.556  # J134
$ws.Cells.Item(134, 11).Value = 13572.522  # K134
$ws.Cells.Item(134, 12).Value = 38833.66800000001  # L134
$ws.Cells.Item(134, 13).Value = -11037.522  # M134
$ws.Cells.Item(134, 14).Value = -43903.66800000001  # N134

$ws.Cells.Item(135, 8).Value = 99941.766  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 99941.766  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 99941.766  # L135
$ws.Cells.Item(135, 14).Value = -110081.766  # N135

$ws.Cells.Item(138, 8).Value = 99995.336  # H138
$ws.Cells.Item(138, 9).Value = 0  # I138
$ws.Cells.Item(138, 10).Value = 99995.336  # J138
$ws.Cells.Item(138, 11).Value = 0  # K138
$ws.Cells.Item(138, 12).Value = 99995.336  # L138
$ws.Cells.Item(138, 14).Value = -110275.336  # N138

$ws.Cells.Item(139, 8).Value = 59998.285  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 59998.285  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 59998.285  # L139
$ws.Cells.Item(139, 14).Value = -70278.285  # N139

$ws.Cells.Item(140, 8).Value = 90999.1  # H140
$ws.Cells.Item(140, 9).Value = 0  # I140
$ws.Cells.Item(140, 10).Value = 90999.1  # J140
$ws.Cells.Item(140, 11).Value = 0  # K140
$ws.Cells.Item(140, 12).Value = 90999.1  # L140
$ws.Cells.Item(140, 14).Value = -101359.1  # N140

$ws.Cells.Item(141, 8).Value = 59900.086  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 59900.086  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 59900.086  # L141
$ws.Cells.Item(141, 14).Value = -70260.08600000001  # N141

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(31, 8).Value = 2235.6897  # H31
$ws.Cells.Item(31, 9).Value = 1631.4584  # I31
$ws.Cells.Item(31, 10).Value = 5136  # J31
$ws.Cells.Item(31, 11).Value = 1631.4584  # K31
$ws.Cells.Item(31, 12).Value = 5136  # L31
$ws.Cells.Item(31, 13).Value = -1336.4584  # M31
$ws.Cells.Item(31, 14).Value = -5726  # N31

$ws.Cells.Item(34, 8).Value = 2235.6897  # H34
$ws.Cells.Item(34, 9).Value = 1631.4584  # I34
$ws.Cells.Item(34, 10).Value = 5136  # J34
$ws.Cells.Item(34, 11).Value = 1631.4584  # K34
$ws.Cells.Item(34, 12).Value = 5136  # L34
$ws.Cells.Item(34, 13).Value = -1429.4584  # M34
$ws.Cells.Item(34, 14).Value = -5540  # N34

$ws.Cells.Item(105, 8).Value = 2744.111  # H105
$ws.Cells.Item(105, 9).Value = 2478  # I105
$ws.Cells.Item(105, 10).Value = 3276.3333  # J105
$ws.Cells.Item(105, 11).Value = 2478  # K105
$ws.Cells.Item(105, 12).Value = 3276.3333  # L105
$ws.Cells.Item(105, 13).Value = -731  # M105
$ws.Cells.Item(105, 14).Value = -6770.3333  # N105

$ws.Cells.Item(107, 8).Value = 572.8788  # H107
$ws.Cells.Item(107, 9).Value = 497.03125  # I107
$ws.Cells.Item(107, 10).Value = 3000  # J107
$ws.Cells.Item(107, 11).Value = 497.03125  # K107
$ws.Cells.Item(107, 12).Value = 3000  # L107
$ws.Cells.Item(107, 13).Value = 1422.96875  # M107
$ws.Cells.Item(107, 14).Value = -6840  # N107

$ws.Cells.Item(122, 8).Value = 1940.7273  # H122
$ws.Cells.Item(122, 9).Value = 2018.125  # I122
$ws.Cells.Item(122, 10).Value = 1734.3334  # J122
$ws.Cells.Item(122, 11).Value = 6054.375  # K122
$ws.Cells.Item(122, 12).Value = 5203.0002  # L122
$ws.Cells.Item(122, 13).Value = -3604.375  # M122
$ws.Cells.Item(122, 14).Value = -10103.0002  # N122

$ws.Cells.Item(132, 8).Value = 25198.732  # H132
$ws.Cells.Item(132, 9).Value = 18602.229  # I132
$ws.Cells.Item(132, 10).Value = 34433.84  # J132
$ws.Cells.Item(132, 11).Value = 55806.687  # K132
$ws.Cells.Item(132, 12).Value = 103301.52  # L132
$ws.Cells.Item(132, 13).Value = -53276.687  # M132
$ws.Cells.Item(132, 14).Value = -108361.52  # N132

$ws.Cells.Item(134, 8).Value = 6293.394  # H134
$ws.Cells.Item(134, 9).Value = 5396.56  # I134
$ws.Cells.Item(134, 10).Value = 9096  # J134
$ws.Cells.Item(134, 11).Value = 16189.68  # K134
$ws.Cells.Item(134, 12).Value = 27288  # L134
$ws.Cells.Item(134, 13).Value = -13654.68  # M134
$ws.Cells.Item(134, 14).Value = -32358  # N134

$ws.Cells.Item(141, 8).Value = 85998.71  # H141
$ws.Cells.Item(141, 9).Value = 105000  # I141
$ws.Cells.Item(141, 10).Value = 84271.32  # J141
$ws.Cells.Item(141, 11).Value = 105000  # K141
$ws.Cells.Item(141, 12).Value = 84271.32  # L141
$ws.Cells.Item(141, 13).Value = -99820  # M141
$ws.Cells.Item(141, 14).Value = -94631.32  # N141

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(5, 8).Value = 2083.6843  # H5
$ws.Cells.Item(5, 9).Value = 1483.3334  # I5
$ws.Cells.Item(5, 10).Value = 2196.25  # J5
$ws.Cells.Item(5, 11).Value = 4450.0002  # K5
$ws.Cells.Item(5, 12).Value = 6588.75  # L5
$ws.Cells.Item(5, 13).Value = -4338.0002  # M5
$ws.Cells.Item(5, 14).Value = -6812.75  # N5

$ws.Cells.Item(47, 8).Value = 306.4  # H47
$ws.Cells.Item(47, 9).Value = 306.4  # I47
$ws.Cells.Item(47, 10).Value = 0  # J47
$ws.Cells.Item(47, 11).Value = 919.1999999999999  # K47
$ws.Cells.Item(47, 12).Value = 0  # L47
$ws.Cells.Item(47, 13).Value = -488.1999999999999  # M47

$ws.Cells.Item(68, 8).Value = 1999.6666  # H68
$ws.Cells.Item(68, 9).Value = 0  # I68
$ws.Cells.Item(68, 10).Value = 1999.6666  # J68
$ws.Cells.Item(68, 11).Value = 0  # K68
$ws.Cells.Item(68, 12).Value = 5998.9998  # L68
$ws.Cells.Item(68, 14).Value = -7620.9998  # N68

$ws.Cells.Item(71, 8).Value = 1999.6666  # H71
$ws.Cells.Item(71, 9).Value = 0  # I71
$ws.Cells.Item(71, 10).Value = 1999.6666  # J71
$ws.Cells.Item(71, 11).Value = 0  # K71
$ws.Cells.Item(71, 12).Value = 17996.9994  # L71
$ws.Cells.Item(71, 14).Value = -26108.9994  # N71

$ws.Cells.Item(134, 8).Value = 2241.6667  # H134
$ws.Cells.Item(134, 9).Value = 2241.6667  # I134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 11).Value = 6725.000100000001  # K134
$ws.Cells.Item(134, 12).Value = 0  # L134
$ws.Cells.Item(134, 13).Value = -1655.000100000001  # M134

$ws.Cells.Item(135, 8).Value = 2083.6843  # H135
$ws.Cells.Item(135, 9).Value = 1483.3334  # I135
$ws.Cells.Item(135, 10).Value = 2196.25  # J135
$ws.Cells.Item(135, 11).Value = 13350.0006  # K135
$ws.Cells.Item(135, 12).Value = 19766.25  # L135
$ws.Cells.Item(135, 13).Value = -10815.0006  # M135
$ws.Cells.Item(135, 14).Value = -24836.25  # N135

$ws.Cells.Item(137, 8).Value = 2714.2354  # H137
$ws.Cells.Item(137, 9).Value = 2496.6667  # I137
$ws.Cells.Item(137, 10).Value = 2832.9092  # J137
$ws.Cells.Item(137, 11).Value = 7490.000100000001  # K137
$ws.Cells.Item(137, 12).Value = 8498.7276  # L137
$ws.Cells.Item(137, 13).Value = -2390.000100000001  # M137
$ws.Cells.Item(137, 14).Value = -18698.7276  # N137

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Cells.Item(122, 8).Value = 6666  # H122
$ws.Cells.Item(122, 9).Value = 4444  # I122
$ws.Cells.Item(122, 10).Value = 8888  # J122
$ws.Cells.Item(122, 11).Value = 13332  # K122
$ws.Cells.Item(122, 12).Value = 26664  # L122
$ws.Cells.Item(122, 13).Value = -10882  # M122
$ws.Cells.Item(122, 14).Value = -31564  # N122

$ws.Cells.Item(126, 8).Value = 4898.25  # H126
$ws.Cells.Item(126, 9).Value = 4362.2  # I126
$ws.Cells.Item(126, 10).Value = 5281.143  # J126
$ws.Cells.Item(126, 11).Value = 13086.6  # K126
$ws.Cells.Item(126, 12).Value = 15843.429  # L126
$ws.Cells.Item(126, 13).Value = -10616.6  # M126
$ws.Cells.Item(126, 14).Value = -20783.429  # N126

$ws.Cells.Item(132, 8).Value = 3148.6667  # H132
$ws.Cells.Item(132, 9).Value = 2696.647  # I132
$ws.Cells.Item(132, 10).Value = 5069.75  # J132
$ws.Cells.Item(132, 11).Value = 8089.941  # K132
$ws.Cells.Item(132, 12).Value = 15209.25  # L132
$ws.Cells.Item(132, 13).Value = -5559.941  # M132
$ws.Cells.Item(132, 14).Value = -20269.25  # N132

$ws.Cells.Item(135, 8).Value = 49083  # H135
$ws.Cells.Item(135, 9).Value = 0  # I135
$ws.Cells.Item(135, 10).Value = 49083  # J135
$ws.Cells.Item(135, 11).Value = 0  # K135
$ws.Cells.Item(135, 12).Value = 49083  # L135
$ws.Cells.Item(135, 14).Value = -59223  # N135

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(22, 8).Value = 2489.4  # H22
$ws.Cells.Item(22, 9).Value = 2298.9092  # I22
$ws.Cells.Item(22, 10).Value = 2722.2222  # J22
$ws.Cells.Item(22, 11).Value = 2298.9092  # K22
$ws.Cells.Item(22, 12).Value = 2722.2222  # L22
$ws.Cells.Item(22, 13).Value = -2003.9092  # M22
$ws.Cells.Item(22, 14).Value = -3312.2222  # N22

$ws.Cells.Item(27, 8).Value = 2489.4  # H27
$ws.Cells.Item(27, 9).Value = 2298.9092  # I27
$ws.Cells.Item(27, 10).Value = 2722.2222  # J27
$ws.Cells.Item(27, 11).Value = 2298.9092  # K27
$ws.Cells.Item(27, 12).Value = 2722.2222  # L27
$ws.Cells.Item(27, 13).Value = -2191.9092  # M27
$ws.Cells.Item(27, 14).Value = -2936.2222  # N27

$ws.Cells.Item(40, 8).Value = 2164.889  # H40
$ws.Cells.Item(40, 9).Value = 2056.9412  # I40
$ws.Cells.Item(40, 10).Value = 4000  # J40
$ws.Cells.Item(40, 11).Value = 2056.9412  # K40
$ws.Cells.Item(40, 12).Value = 4000  # L40
$ws.Cells.Item(40, 13).Value = -1920.9412  # M40
$ws.Cells.Item(40, 14).Value = -4272  # N40

$ws.Cells.Item(82, 8).Value = 1410  # H82
$ws.Cells.Item(82, 9).Value = 1524.9166  # I82
$ws.Cells.Item(82, 10).Value = 950.3333  # J82
$ws.Cells.Item(82, 11).Value = 1524.9166  # K82
$ws.Cells.Item(82, 12).Value = 950.3333  # L82
$ws.Cells.Item(82, 13).Value = -1163.9166  # M82
$ws.Cells.Item(82, 14).Value = -1672.3333  # N82

$ws.Cells.Item(85, 8).Value = 1410  # H85
$ws.Cells.Item(85, 9).Value = 1524.9166  # I85
$ws.Cells.Item(85, 10).Value = 950.3333  # J85
$ws.Cells.Item(85, 11).Value = 1524.9166  # K85
$ws.Cells.Item(85, 12).Value = 950.3333  # L85
$ws.Cells.Item(85, 13).Value = -276.9166  # M85
$ws.Cells.Item(85, 14).Value = -3446.3333  # N85

$ws.Cells.Item(132, 8).Value = 7939502.5  # H132
$ws.Cells.Item(132, 9).Value = 10103894  # I132
$ws.Cells.Item(132, 10).Value = 3402.6667  # J132
$ws.Cells.Item(132, 11).Value = 30311682  # K132
$ws.Cells.Item(132, 12).Value = 10208.0001  # L132
$ws.Cells.Item(132, 13).Value = -30309152  # M132
$ws.Cells.Item(132, 14).Value = -15268.0001  # N132

$ws.Cells.Item(136, 8).Value = 4833076  # H136
$ws.Cells.Item(136, 9).Value = 5557665.5  # I136
$ws.Cells.Item(136, 10).Value = 2478.3333  # J136
$ws.Cells.Item(136, 11).Value = 16672996.5  # K136
$ws.Cells.Item(136, 12).Value = 7434.999899999999  # L136
$ws.Cells.Item(136, 13).Value = -16670446.5  # M136
$ws.Cells.Item(136, 14).Value = -12534.9999  # N136

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(132, 8).Value = 14734.492  # H132
$ws.Cells.Item(132, 9).Value = 9259.072  # I132
$ws.Cells.Item(132, 10).Value = 39830.168  # J132
$ws.Cells.Item(132, 11).Value = 27777.216  # K132
$ws.Cells.Item(132, 12).Value = 119490.504  # L132
$ws.Cells.Item(132, 13).Value = -25247.216  # M132
$ws.Cells.Item(132, 14).Value = -124550.504  # N132
